$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutWs = $wb.Worksheets.Item("About")
$dataWs  = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the two cells on the "About" sheet that embed the build timestamp.
$a2 = $aboutWs.Range("A2")
$a2.Value = $a2.Value().Replace($oldStamp, $newStamp)

$a6 = $aboutWs.Range("A6")
$a6.Value = $a6.Value().Replace($oldStamp, $newStamp)

# Update the build_version column (S2:S147) on the data sheet.
$lastRow = $dataWs.Cells.Item($dataWs.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataWs.Cells.Item($r, 19)
    $val = $cell.Value()
    if ($val -and $val.ToString().Contains($oldStamp)) {
        $cell.Value = $val.ToString().Replace($oldStamp, $newStamp)
    }
}
